# Renovacion informe tecnico - agrega el campo "ultimo_convenio" a los
# placeholders de comidas aprobadas (Desayuno/Almuerzo/Merienda/Cena x
# dia de la semana) en la primera tabla de "prestaciones financiadas",
# y actualiza los campos de pagina en el pie (el documento paso a tener
# una pagina mas).

$d = $word.ActiveDocument

# Solo la primera tabla semanal ("financiaron las siguientes
# prestaciones") debe cambiar; la segunda tabla semanal ("cantidades
# aprobadas" en la Evaluacion Tecnica) permanece igual.
$tabla = $d.Tables.Item(4)

$dias = "lunes", "martes", "miercoles", "jueves", "viernes", "sabado", "domingo"
$comidas = "desayuno", "almuerzo", "merienda", "cena"

$primera = $true
foreach ($dia in $dias) {
    foreach ($comida in $comidas) {
        $viejo = "informe.aprobadas_${comida}_${dia}"
        if ($primera) {
            $nuevo = "informe.aprobadas_ultimo_convenio_${comida}_${dia}"
            $primera = $false
        } else {
            $nuevo = "informe.aprobadas_ ultimo_convenio_${comida}_${dia}"
        }

        $rng = $tabla.Range
        $rng.Find.Execute($viejo, $true, $false, $false, $false, $false, $true, 1, $false, $nuevo, 2)
    }
}

# El pie de pagina cachea el numero de pagina actual y el total de
# paginas ("Pagina {PAGE} | {NUMPAGES}"); el documento ahora ocupa 7
# paginas en lugar de 6.
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("6", $true, $false, $false, $false, $false, $true, 1, $false, "7", 2)
